$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data set. It belongs right
# above the current row 385, so insert a new row there (this pushes the
# old rows 385-425 down to 386-426, growing the used range to A1:R426).
$ws.Rows.Item(385).Insert()

# Populate the newly inserted row 385 with the new record's values.
$ws.Range("A385").Value = 10
$ws.Range("B385").Value = "Vega Modelo de Temuco"
$ws.Range("C385").Value = "La Araucanía"
$ws.Range("D385").Value = 44918
$ws.Range("E385").Value = 9
$ws.Range("F385").Value = 100112044
$ws.Range("G385").Value = "Perejil"
$ws.Range("H385").Value = "Sin especificar"
$ws.Range("I385").Value = "Primera"
$ws.Range("J385").Value = 65
$ws.Range("K385").Value = 4000
$ws.Range("L385").Value = 4000
$ws.Range("M385").Value = 4000
$ws.Range("N385").Value = "$/docena de atados (3 kilos)"
$ws.Range("O385").Value = "Provincia de Cautín"
$ws.Range("P385").Value = 1333
$ws.Range("Q385").Value = 3
$ws.Range("R385").Value = "Hortaliza"
